# BMRESP BAU Min Req EV Sales Perc.xlsx
# Commit: "Set required EV sales to 0"
#
# The required minimum EV sales percentage trajectory is removed and
# replaced with flat 0% for all years, and the "About" sheet's
# Source/Notes section is rewritten to reflect that no source is
# needed and there is no minimum required EV sales percentage.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Remove the old citation rows (4-9): DRAFT report title, date,
# page number, the hyperlink to the CARB report, and the UCS link.
$about.Rows("4:14").Delete()

# New "Source:" value - no source needed any more.
$about.Range("B3").Value = "none needed"

# Rebuild the Notes section with the new explanation.
$about.Range("A8").Value = "Notes"
$about.Range("A8").Font.Bold = $true
$about.Range("A9").Value = "There is no minimum required EV sales percentage"
$about.Range("A10").Value = "in the United States for any vehicle type."

# Drop the now-unused hyperlink to the old CARB source document.
$about.Hyperlinks.Delete()

# ---------------------------------------------------------------
# Sheet "BMRESP-passenger"
# ---------------------------------------------------------------
$passenger = $wb.Worksheets.Item("BMRESP-passenger")

# Add a header label above the data table describing the units.
$passenger.Range("A1").Value = "Sales Percentage (dimensionless)"
$passenger.Range("A1").Font.Bold = $true
$passenger.Range("A1").WrapText = $true

# LDVs (row 2) no longer ramps up to an 8% minimum requirement -
# set the whole required-sales trajectory (2018-2050) to 0%.
$passenger.Range("D2:AJ2").ClearFormats()
$passenger.Range("D2:AJ2").Value = 0

# ---------------------------------------------------------------
# Sheet "BMRESP-freight"
# ---------------------------------------------------------------
$freight = $wb.Worksheets.Item("BMRESP-freight")

# Same header label as the passenger sheet; the underlying data on
# this sheet was already entirely 0, so no other edits are needed.
$freight.Range("A1").Value = "Sales Percentage (dimensionless)"
$freight.Range("A1").Font.Bold = $true
$freight.Range("A1").WrapText = $true
